# Applies the cryptos-list refresh described by the commit diff.
# Numeric-looking text values (prices) are forced to Text format first so
# Excel stores them as strings (matching the source inlineStr cells) rather
# than auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.271.17"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.126.39"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.07"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "613.85"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.11"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.846"
$ws.Range("E10").Value = "  +14.24%  "
$ws.Range("D11").Value = "3.123.86"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.199"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.39"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").Value = "93.063.14"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.44"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "3.708.74"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "3.125.63"
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.81"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("E21").Value = "  +4.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000202"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "442.91"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.19"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.73"
$ws.Range("E25").Value = "  -0.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.00"
$ws.Range("E26").Value = "  +11.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "86.15"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").Value = "3.293.62"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.182"
$ws.Range("E30").Value = "  +9.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.239"
$ws.Range("E31").Value = "  +3.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("E32").Value = "  -12.77%  "
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.27"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.09"
$ws.Range("E35").Value = "  +6.23%  "
$ws.Range("E36").Value = "  -9.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.03"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.97"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.449"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.99"
$ws.Range("E42").Value = "  +8.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "475.70"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.67"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.701"
$ws.Range("E47").Value = "  +0.98%  "
$ws.Range("E48").Value = "  -1.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.34"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.46"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.06"
$ws.Range("E51").Value = "  +0.02%  "
